$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.479.68"
$ws.Range("E2").Value = "  +2.16%  "

$ws.Range("D3").Value = "2.643.26"
$ws.Range("E3").Value = "  +0.94%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.43%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.177"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.39%  "

$ws.Range("D10").Value = "2.642.06"
$ws.Range("E10").Value = "  +0.93%  "

$ws.Range("E11").Value = "  +1.64%  "

$ws.Range("E12").Value = "  +3.59%  "

$ws.Range("E13").Value = "  +0.27%  "

$ws.Range("E14").Value = "  +5.06%  "

$ws.Range("D15").Value = "3.127.08"
$ws.Range("E15").Value = "  +1.25%  "

$ws.Range("D16").Value = "72.371.91"
$ws.Range("E16").Value = "  +1.95%  "

$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("D18").Value = "2.640.86"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "380.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("E22").Value = "  +0.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.30%  "

$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.17%  "

$ws.Range("E27").Value = "  +4.00%  "

$ws.Range("D28").Value = "2.779.45"
$ws.Range("E28").Value = "  +1.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").Value = "0.0₃0958"
$ws.Range("E30").Value = "  +1.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "524.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.59%  "

$ws.Range("E33").Value = "  -0.50%  "

$ws.Range("E34").Value = "  -0.37%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.46%  "

$ws.Range("E37").Value = "  +1.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.112"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.77%  "

$ws.Range("E39").Value = "  +0.87%  "

$ws.Range("E40").Value = "  +1.89%  "

$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("E42").Value = "  +1.47%  "

$ws.Range("E43").Value = "  +2.03%  "

$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("E45").Value = "  +0.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.543"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.01%  "

$ws.Range("E50").Value = "  +2.76%  "

$ws.Range("E51").Value = "  -3.99%  "
